$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue 2 4 '91.101.28'
Set-TextValue 2 5 '  +0.62%  '
Set-TextValue 3 4 '3.157.54'
Set-TextValue 3 5 '  +1.50%  '
Set-TextValue 4 5 '  +0.12%  '
Set-TextValue 5 4 '238.40'
Set-TextValue 5 5 '  +1.71%  '
Set-TextValue 6 4 '621.00'
Set-TextValue 6 5 '  -0.39%  '
Set-TextValue 7 4 '1.13'
Set-TextValue 7 5 '  +4.93%  '
Set-TextValue 8 5 '  +1.18%  '
Set-TextValue 9 4 '0.999'
Set-TextValue 9 5 '  -0.07%  '
Set-TextValue 10 4 '3.157.98'
Set-TextValue 10 5 '  +1.67%  '
Set-TextValue 11 5 '  +0.84%  '
Set-TextValue 12 4 '0.205'
Set-TextValue 13 5 '  -1.51%  '
Set-TextValue 14 4 '35.39'
Set-TextValue 14 5 '  -1.67%  '
Set-TextValue 15 4 '5.51'
Set-TextValue 15 5 '  +0.60%  '
Set-TextValue 16 4 '91.331.11'
Set-TextValue 16 5 '  +1.20%  '
Set-TextValue 17 4 '3.748.08'
Set-TextValue 17 5 '  +2.36%  '
Set-TextValue 18 4 '3.177.77'
Set-TextValue 18 5 '  +2.75%  '
Set-TextValue 19 5 '  -5.06%  '
Set-TextValue 20 4 '15.28'
Set-TextValue 20 5 '  +8.66%  '
Set-TextValue 21 4 '5.88'
Set-TextValue 21 5 '  +5.41%  '
Set-TextValue 22 4 '0.0000209'
Set-TextValue 22 5 '  -5.96%  '
Set-TextValue 23 4 '442.60'
Set-TextValue 23 5 '  +0.92%  '
Set-TextValue 24 4 '9.19'
Set-TextValue 24 5 '  +2.53%  '
Set-TextValue 25 4 '6.07'
Set-TextValue 25 5 '  +2.50%  '
Set-TextValue 26 4 '88.82'
Set-TextValue 26 5 '  +0.07%  '
Set-TextValue 27 5 '  -0.04%  '
Set-TextValue 28 5 '  +3.21%  '
Set-TextValue 29 4 '0.997'
Set-TextValue 29 5 '  -0.33%  '
Set-TextValue 30 5 '  +47.06%  '
Set-TextValue 31 5 '  +6.48%  '
Set-TextValue 32 5 '  +16.83%  '
Set-TextValue 33 4 '9.36'
Set-TextValue 33 5 '  -0.04%  '
Set-TextValue 34 4 '0.166'
Set-TextValue 34 5 '  +8.92%  '
Set-TextValue 35 4 '7.73'
Set-TextValue 35 5 '  +5.40%  '
Set-TextValue 36 4 '26.39'
Set-TextValue 36 5 '  +1.86%  '
Set-TextValue 37 4 '0.888'
Set-TextValue 37 5 '  -11.17%  '
Set-TextValue 38 4 '509.69'
Set-TextValue 38 5 '  +0.99%  '
Set-TextValue 39 2 'PancakeSwap'
Set-TextValue 39 3 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 39 4 '1.94'
Set-TextValue 39 5 '  +1.63%  '
Set-TextValue 40 2 'Fetch.AI'
Set-TextValue 40 3 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 40 4 '1.36'
Set-TextValue 40 5 '  +5.41%  '
Set-TextValue 41 2 'MantraDAO'
Set-TextValue 41 3 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-TextValue 41 4 '3.88'
Set-TextValue 41 5 '  +12.74%  '
Set-TextValue 42 2 'PolygonEcosystemToken'
Set-TextValue 42 3 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 42 4 '0.453'
Set-TextValue 42 5 '  +11.11%  '
Set-TextValue 43 5 '  -11.34%  '
Set-TextValue 44 4 '22.12'
Set-TextValue 44 5 '  -0.23%  '
Set-TextValue 46 4 '0.717'
Set-TextValue 46 5 '  +3.62%  '
Set-TextValue 47 4 '156.62'
Set-TextValue 47 5 '  +3.60%  '
Set-TextValue 48 5 '  +0.73%  '
Set-TextValue 49 5 '  +3.83%  '
Set-TextValue 50 2 'Filecoin'
Set-TextValue 50 3 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 50 4 '4.45'
Set-TextValue 50 5 '  -0.12%  '
Set-TextValue 51 2 'OKB'
Set-TextValue 51 3 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 51 4 '44.02'
Set-TextValue 51 5 '  -1.28%  '
